$wb = $excel.ActiveWorkbook

# --- SYGBSC: shift starting year from 2020 to 2021, referencing Summary column C instead of B ---
$wsSygbsc = $wb.Worksheets.Item("SYGBSC")
$wsSygbsc.Range("B1").Value = 2021
$wsSygbsc.Range("B2").Formula = "=Summary!C6+Summary!C13"

# --- BGBSC: drop the (now-duplicated) first projection year column (2021) ---
$wsBgbsc = $wb.Worksheets.Item("BGBSC")
$wsBgbsc.Range("B:B").EntireColumn.Delete()

# --- Cursor / selection bookkeeping to match the saved view state ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Select()

$wsBgbsc.Activate()
$wsBgbsc.Range("E21").Select()
